# Updated cryptos list on Tue Aug  6 19:41:20 UTC 2024 with GitHub Actions
#
# Refresh of the Price (D) / Volume(1h) (E) columns for most rows, plus a
# handful of rows whose rank changed (their Coin name + Link swapped with
# the neighbouring row): 42/43 (Mantle <-> FirstDigitalUSD),
# 46/47 (Stellar <-> Bittensor), 50/51 (EnergySwap <-> Maker).
#
# Price/Volume cells are stored as text in the workbook (e.g. "56.815.94",
# "1.00", "  +6.83%  "), so for the handful of new Price values that look
# like plain numbers ("490.25", "1.00", ...) we force the cell to Text
# format before assigning, otherwise Excel would silently coerce them to
# numeric values (losing formatting like trailing zeros) - then restore
# the cell's style so no stray formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '56.815.94'
$ws.Range("E2").Value = '  +6.83%  '
$ws.Range("D3").Value = '2.484.55'
$ws.Range("E3").Value = '  +4.15%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '490.25'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +7.21%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.54'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +14.47%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.997'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.516'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +7.57%  '
$ws.Range("D9").Value = '2.502.79'
$ws.Range("E9").Value = '  +4.19%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.79'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +10.40%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0982'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.38%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.332'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +7.38%  '
$ws.Range("E13").Value = '  +1.89%  '
$ws.Range("D14").Value = '2.916.29'
$ws.Range("E14").Value = '  +4.18%  '
$ws.Range("D15").Value = '56.586.40'
$ws.Range("E15").Value = '  +6.27%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.25'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +9.27%  '
$ws.Range("E17").Value = '  +5.56%  '
$ws.Range("D18").Value = '2.502.21'
$ws.Range("E18").Value = '  +4.66%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.56'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +10.82%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.12'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +9.33%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '319.67'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +5.10%  '
$ws.Range("E22").Value = '  +0.28%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.83'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +10.34%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '58.93'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +6.26%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.411'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +8.12%  '
$ws.Range("E26").Value = '  -0.84%  '
$ws.Range("E27").Value = '  +8.45%  '
$ws.Range("D28").Value = '2.592.14'
$ws.Range("E28").Value = '  +5.01%  '
$ws.Range("E29").Value = '  +9.48%  '
$ws.Range("E30").Value = '  +11.00%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.31%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '149.28'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.13%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.24'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.10%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.23'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +5.86%  '
$ws.Range("E36").Value = '  +9.51%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.73'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +7.14%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.860'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +8.85%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '34.24'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.63%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.51'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +8.33%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0561'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +8.02%  '
$ws.Range("B42").Value = 'FirstDigitalUSD'
$ws.Range("C42").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.997'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.10%  '
$ws.Range("B43").Value = 'Mantle'
$ws.Range("C43").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.611'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.88%  '
$ws.Range("E44").Value = '  +8.95%  '
$ws.Range("E45").Value = '  +15.50%  '
$ws.Range("B46").Value = 'Bittensor'
$ws.Range("C46").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '259.41'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +19.09%  '
$ws.Range("B47").Value = 'Stellar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0921'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +7.37%  '
$ws.Range("E48").Value = '  +6.22%  '
$ws.Range("E49").Value = '  +0.73%  '
$ws.Range("B50").Value = 'Maker'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D50").Value = '1.895.29'
$ws.Range("E50").Value = '  -2.19%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '17.62'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +8.00%  '
